$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AF4").Value = 0.717
$ws.Range("AF5").Value = 0.957
$ws.Range("AF6").Value = 0.82
$ws.Range("AF7").Value = 0.897
$ws.Range("AF8").Value = 0.887
$ws.Range("AF9").Value = 0.783
$ws.Range("AF10").Value = 0.957
$ws.Range("AF11").Value = 0.957
$ws.Range("AF12").Value = 1.227
$ws.Range("AF13").Value = 1.565
